$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 - same look (bold + border) as the other header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# New data columns I (I0) and J (IF) for rows 2-13
$data = @(
    @(9, 9),
    @(3, 3),
    @(4, 6),
    @(5, 6),
    @(8, 9),
    @(3, 5),
    @(1, 4),
    @(1, 4),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(6, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
